# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets. Each worksheet gets the same set of row updates.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 1922
    7  = 1602
    9  = 628
    13 = 92
    18 = 125
    19 = 3706
    20 = 5
    21 = 7
    23 = 338
    24 = 593
    25 = 362
    27 = 29
    28 = 1502
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
